# Update "想去人数" (F column) counts across the sheets to reflect the
# newly generated data snapshot (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$s1 = @{
    2  = 202
    3  = 113
    5  = 984
    6  = 5452
    7  = 481
    8  = 675
    9  = 942
    11 = 77
    13 = 584
    14 = 25
    17 = 1818
    18 = 1462
    19 = 883
    22 = 328
    23 = 536
    24 = 144
    25 = 1053
    28 = 2822
    33 = 33
    34 = 363
    37 = 11
    39 = 286
    40 = 689
    41 = 86
    43 = 54
    44 = 66
}
foreach ($row in $s1.Keys) {
    $ws1.Range("F$row").Value = $s1[$row]
}

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$s2 = @{
    6  = 127
    10 = 11
}
foreach ($row in $s2.Keys) {
    $ws2.Range("F$row").Value = $s2[$row]
}

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$s4 = @{
    3  = 202
    4  = 113
    5  = 984
    7  = 5452
    8  = 481
    9  = 675
    12 = 942
    15 = 127
    16 = 77
    18 = 584
    19 = 25
    23 = 1818
    24 = 1462
    25 = 884
    27 = 328
    29 = 536
    30 = 144
    31 = 1053
    32 = 2822
    37 = 33
    38 = 363
    41 = 11
    42 = 286
    43 = 689
    44 = 86
    45 = 54
    46 = 66
    48 = 11
}
foreach ($row in $s4.Keys) {
    $ws4.Range("F$row").Value = $s4[$row]
}

$wb.Save()
